$wb = $excel.ActiveWorkbook

# --- Reorder worksheets: review_info first, hotel_info second ---
$hotelWs = $wb.Worksheets.Item("hotel_info")
$reviewWs = $wb.Worksheets.Item("review_info")
$hotelWs.Move($null, $reviewWs)

# --- hotel_info: insert a new "State" column between Hotel_Name and City ---
$ws = $wb.Worksheets.Item("hotel_info")
$ws.Columns.Item(3).Insert()
$ws.Cells.Item(1,3).Value = "State"
$ws.Cells.Item(2,3).Value = "Louisiana"
